# Apply updated dSF (column F) values per the "repull data, push all data,
# mean calculation" commit. Only column F changes; column E (dS0) stays
# exactly as-is in every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    5  = 1
    7  = -1
    8  = 2
    9  = -1
    11 = -2
    13 = 4
    14 = -2
    17 = 1
    19 = -1
    22 = -2
    26 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
